# "Activité 3 DdF TPA3 du 16/10" — add the 2023-10-16 journal entry (row 22)
# to the "2023-2024" sheet of the Journal de bord.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 22 was a blank templated row (style only, no content). Pull the
# correctly-formatted styles from neighbouring filled-in cells so the new
# row ends up with the same cellXfs (date format on A22, wrap+vcenter on
# I22) as the rest of the table instead of Excel minting new style ids.
$ws.Range("A21").Copy() | Out-Null
$ws.Range("A22").PasteSpecial(-4122) | Out-Null   # xlPasteFormats (date format)
$ws.Range("I21").Copy() | Out-Null
$ws.Range("I22").PasteSpecial(-4122) | Out-Null   # xlPasteFormats (wrap + vcenter)
$excel.CutCopyMode = 0

# Values for the new entry.
$ws.Range("A22").Value = 45215                    # 2023-10-16
$ws.Range("B22").Value = "MPAL"
$ws.Range("C22").Value = "TP"
$ws.Range("F22").Value = "x"
$ws.Range("G22").Value = "Bilan et synthèse autour de l'utilisation des Business Rules vs Tests d'Acceptation.`nActivité 3 : Définition de Fini"
$ws.Range("I22").Value = "Projet Web de la SAE pas encore lancé. Travaillé sur le projet Java de la SAE.`nQuestion : ""Les tests unitaires passent"" et ""Les tests d'acceptation passent"" vs ""Au moins 80% de couverture fonctionnelle""?`nJe ne vois pas l'intérêt de dire que tous 100% des tests unitaires écrits passent : si certains ne passent pas, on les supprime et on arrive à nouveau à 100%...`nPour 100% des tests d'acceptation : ne pas les respecter reviendrait à ne pas terminer l'implémentation de l'US non ? ce qui me parait étrange...`nPar ailleurs, pourrait-tu m'aider à clarifier ce qui signifie le niveau de couverture par les tests ? est-ce équivalent à la couverture fonctionnelle ?`n"

# Reflect where the author ended up after filling the row (I21 -> I23).
$ws.Range("I23").Select() | Out-Null
